$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.021.09'
$ws.Range('E2').Value = '  -2.32%  '
$ws.Range('D3').Value = '2.428.69'
$ws.Range('E3').Value = '  -1.12%  '
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').Value = '''572.54'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.34%  '
$ws.Range('D6').Value = '''140.65'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -2.25%  '
$ws.Range('E7').Value = '  +0.17%  '
$ws.Range('D8').Value = '''0.528'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -1.04%  '
$ws.Range('D9').Value = '2.415.35'
$ws.Range('E9').Value = '  -1.52%  '
$ws.Range('E10').Value = '  -0.86%  '
$ws.Range('E11').Value = '  -0.09%  '
$ws.Range('D12').Value = '''5.12'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -2.05%  '
$ws.Range('E13').Value = '  -1.98%  '
$ws.Range('D14').Value = '''26.10'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -1.51%  '
$ws.Range('D15').Value = '''0.0000172'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -4.20%  '
$ws.Range('D16').Value = '2.826.54'
$ws.Range('E16').Value = '  -2.46%  '
$ws.Range('D17').Value = '60.897.76'
$ws.Range('E17').Value = '  -2.00%  '
$ws.Range('D18').Value = '2.417.30'
$ws.Range('E18').Value = '  -1.46%  '
$ws.Range('D19').Value = '''7.59'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +6.04%  '
$ws.Range('D20').Value = '''10.73'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -1.44%  '
$ws.Range('D21').Value = '''323.71'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -1.57%  '
$ws.Range('D22').Value = '''4.07'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -1.66%  '
$ws.Range('D23').Value = '''6.07'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +1.14%  '
$ws.Range('E24').Value = '  +0.23%  '
$ws.Range('D25').Value = '''1.89'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -4.85%  '
$ws.Range('D26').Value = '''64.95'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -1.24%  '
$ws.Range('D27').Value = '''590.57'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +0.23%  '
$ws.Range('D28').Value = '''8.48'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -8.78%  '
$ws.Range('D29').Value = '2.536.70'
$ws.Range('E29').Value = '  -1.61%  '
$ws.Range('D30').Value = '0.0₃0936'
$ws.Range('E30').Value = '  -3.85%  '
$ws.Range('D31').Value = '''7.96'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -1.52%  '
$ws.Range('E32').Value = '  -5.35%  '
$ws.Range('D33').Value = '''1.86'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -2.38%  '
$ws.Range('E34').Value = '  -1.46%  '
$ws.Range('E35').Value = '  +0.00%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').Value = '''1.42'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -2.12%  '
$ws.Range('B37').Value = 'NEARProtocol'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D37').Value = '''4.67'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -5.66%  '
$ws.Range('D38').Value = '''151.98'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -1.62%  '
$ws.Range('D39').Value = '''0.369'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -2.99%  '
$ws.Range('D40').Value = '''18.34'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -0.75%  '
$ws.Range('D41').Value = '''5.17'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -3.28%  '
$ws.Range('E42').Value = '  +0.10%  '
$ws.Range('D43').Value = '''1.69'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -2.32%  '
$ws.Range('D44').Value = '''41.18'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -4.59%  '
$ws.Range('D45').Value = '''2.37'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -4.64%  '
$ws.Range('D46').Value = '0.0₆0288'
$ws.Range('E46').Value = '  +8.85%  '
$ws.Range('D47').Value = '''142.05'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.92%  '
$ws.Range('D48').Value = '''3.53'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -3.86%  '
$ws.Range('D49').Value = '''0.593'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -2.62%  '
$ws.Range('D50').Value = '''19.74'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -1.00%  '
$ws.Range('D51').Value = '''0.0507'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -3.89%  '
